$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 162.6848836413069
$ws.Range("B3").Value = 49.76516661019377
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 143934.1102782082
$ws.Range("B6").Value = 9707.816839233754
$ws.Range("B7").Value = 8756.876020340693
$ws.Range("B8").Value = 20217.59999999929
$ws.Range("B9").Value = 616.0374081889174
$ws.Range("B10").Value = 45710.66540329998
$ws.Range("B11").Value = 0.07771981492021149
$ws.Range("B12").Value = 0.4250674068788627
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0.9933932212659528
$ws.Range("B15").Value = 0.9364218080283236
